# Refresh the cached regression-table figures in decomposition_main_te.xlsx.
#
# The worksheet's numeric cells are formulas that pull cached text from an
# external workbook reference (=[1]decomposition_main_te!...). The linked
# source file isn't available in this environment, so we can't "recalculate"
# those links the normal way. Instead we push the updated, already-rounded
# figures (as they'd appear after the source refresh) straight into each
# cell, using a self-referential text formula (="value") so the cell keeps
# behaving like a formula/text result (t="str") rather than turning into a
# plain typed-in value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CachedText {
    param(
        [string]$CellRef,
        [string]$NewText
    )
    $escaped = $NewText.Replace('"', '""')
    $ws.Range($CellRef).Formula = '="' + $escaped + '"'
}

# Row 6 (coefficients) / Row 7 (std errors) -- FC regression block
Set-CachedText "B6" "-202.3***"
Set-CachedText "E6" "-0.34"
Set-CachedText "F6" "-77.0**"
Set-CachedText "B7" "(48.2)"
Set-CachedText "E7" "(2.84)"

# Row 8 (coefficients) / Row 9 (std errors) -- Default regression block
Set-CachedText "B8" "-40.1"
Set-CachedText "E8" "-1.43"
Set-CachedText "F8" "-16.5"
Set-CachedText "G8" "-0.023"
Set-CachedText "I8" "-0.0077"
Set-CachedText "E9" "(2.44)"
Set-CachedText "F9" "(33.0)"

# Row 13 -- Control mean
Set-CachedText "B13" "941.1"
Set-CachedText "E13" "5.26"
Set-CachedText "F13" "395.2"
Set-CachedText "G13" "0.43"
